# Add a new "UID" column (K) to the roster sheet, with a sequential
# unique id (3001, 3002, ...) for every data row, mirroring the
# exported SQL "UID" primary key used for the new foreign-key/cascade
# relationship described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cell (K1) the same look as the existing "Lib"
# header in J1 (bordered / centered header style) by copying its
# formatting, then overwrite the copied text with the new header.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Value = "UID"

# Fill K2:K92 with sequential UID values 3001..3091 (one per data row).
for ($i = 2; $i -le 92; $i++) {
    $ws.Cells.Item($i, 11).Value = 3000 + ($i - 1)
}

# Match the author's final selection: column K, rows 2-92 selected.
$ws.Range("K2:K92").Select() | Out-Null
